$d = $word.ActiveDocument

# --- 1) Mark the four inline pictures as "no proofing" (adds <w:noProof/> to
#        the rPr of the runs that host each <w:drawing>). This is what Word
#        stamps onto picture runs once the document is opened/edited in the
#        real app. ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = $true
}

# --- 2) Truth-table correction: "Case 3" (A = 5v, B = 5v) must read 5v,
#        not 0v, for an XNOR gate. Find the paragraph, locate the single
#        "0" character that sits between "op is " and "v", and fix it. ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    if ($text.StartsWith("Case 3") -and $text.Contains("op is 0v")) {
        $rng = $p.Range
        $idx = $rng.Text.IndexOf("op is 0v")
        $charStart = $rng.Start + $idx + 6
        $charRng = $rng.Duplicate
        $charRng.SetRange($charStart, $charStart + 1)
        $charRng.Text = "5"
        break
    }
}

# --- 3) Header correction: the default header reads "NOR Gate" but this
#        document is about the XNOR gate, so prefix it with a bold "X" that
#        matches the existing run formatting. ---
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdr.Range.InsertBefore("X")

$hdr2 = $sec.Headers.Item(1)
$full = $hdr2.Range
$xRng = $full.Duplicate
$xRng.SetRange($full.Start, $full.Start + 1)
$xRng.Bold = $true
$xRng.Font.BoldBi = $true
$xRng.Font.Size = 20
$xRng.Font.SizeBi = 20
$xRng.LanguageID = "en-US"
